$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The laptop model in B3 was renamed to include the "ASUS" brand prefix.
$ws.Range("B3").Value = "ASUS Eee PC 1201NL"

# Update the selected cell/range in the sheet view from C6 to B3.
$ws.Range("B3").Select()
